# Daily attendance processing - rotate "Recorded By" contributor lists
# For each populated cell in column G (rows 2-148), the list of recorders
# (comma-separated) is rotated so the first entry moves to the end.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 148; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $cell.Value = $rotated -join ", "
        }
    }
}
